$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TCDA")

# Row 21
$ws.Range("D21").Value = 400
$ws.Range("E21").Value = "NA"

# Row 83
$ws.Range("D83").Value = 400
$ws.Range("E83").Value = "NA"

# Row 89
$ws.Range("D89").Value = -23100
$ws.Range("E89").Value = "NA"

# Row 91
$ws.Range("D91").Value = -200
$ws.Range("E91").Value = "NA"

# Row 94
$ws.Range("D94").Value = -20300
$ws.Range("E94").Value = "NA"

# Row 100
$ws.Range("D100").Value = 43000
$ws.Range("E100").Value = "NA"

# Row 101
$ws.Range("J101").Value = "NA"

# Row 102
$ws.Range("D102").Value = -400
$ws.Range("E102").Value = "NA"
